$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column K (rows 3-11) into new column L, which extends the table with
# another year of data.
$ws.Range("K3:K11").Copy() | Out-Null
$ws.Range("L3:L11").PasteSpecial() | Out-Null

# The new column represents 2021 (column K was 2020).
$ws.Range("L4").Value = 2021

# Update the active selection to match the post-edit state.
$ws.Range("N2").Select() | Out-Null
